$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    # MatchCase=$true, MatchWholeWord=$false, MatchWildcards=$false, MatchSoundsLike=$false,
    # MatchAllWordForms=$false, Forward=$true, Wrap=1 (wdFindContinue), Format=$false,
    # Replace=1 (wdReplaceOne) -- only touch the first (leftmost) match so that documents with
    # repeated substrings (e.g. a company name mentioned twice) only change the intended spot.
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 1)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
    }
}

# 1. July 2021 start date -> July 26, 2021
Replace-Text " July 2021 - Present (1 year " " July 26, 2021 - Present (1 year "

# 2. Remove spell-check artifact around Javascript (merge runs, same visible text)
Replace-Text " Selenium C# with the C.A.T.S framework to develop automation tests for both the Javascript web portal frontend and API endpoints, seamlessly integrating automation stories into the Azure CLI automation pipeline and generating comprehensive reports on test results " " Selenium C# with the C.A.T.S framework to develop automation tests for both the Javascript web portal frontend and API endpoints, seamlessly integrating automation stories into the Azure CLI automation pipeline and generating comprehensive reports on test results "

# 3. Jan 2021 - June 2021 -> Jan 18, 2021 - June 23, 2021
Replace-Text "Jan 2021 – June 2021 " "Jan 18, 2021 – June 23, 2021 "

# 4. Remove spell-check artifact around NUnit (merge runs, same visible text)
Replace-Text "Developed individual unit tests for each product using NUnit prior to submitting them to QA for quality assurance testing." "Developed individual unit tests for each product using NUnit prior to submitting them to QA for quality assurance testing."

# 5. Nov 2019 - Jan 2021 -> Nov 8, 2019 - Jan 11, 2021
Replace-Text "Nov 2019 – Jan 2021 | " "Nov 8, 2019 – Jan 11, 2021 | "

# 6. Company name change (the commit's headline change)
Replace-Text "Martina Carter Entertainment" "Edutainment Living History"

# 7. Aug 2019 - Dec 2019 -> Aug 01, 2019 - Dec 01, 2019
Replace-Text "Aug 2019 – Dec 2019" "Aug 01, 2019 – Dec 01, 2019"

# 8. Zerma -> Bioenergy Technology Inc (also removes spell-check artifact)
Replace-Text "Zerma " "Bioenergy Technology Inc "

# 9. Remove grammar-check artifact around Bachelor's (merge runs, same visible text)
Replace-Text "Bachelor’s in Computer Science" "Bachelor’s in Computer Science"

Write-Output "All replacements attempted."
